# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 13 de Mayo de 2020 a las 02:05"

# Row 4: Estados Unidos - refreshed totals
$ws.Range("B4").Value = 1408155
$ws.Range("C4").Value = 22321
$ws.Range("D4").Value = 295781
$ws.Range("E4").Value = 1028997
$ws.Range("G4").Value = 1582
$ws.Range("H4").Value = 83377

# Rows 37-38: Japon moves above Austria (table sorted by total cases desc)
$ws.Range("A37").Value = "Japon"
$ws.Range("B37").Value = 15968
$ws.Range("C37").Value = 121
$ws.Range("D37").Value = 8531
$ws.Range("E37").Value = 6780
$ws.Range("F37").Value = 249
$ws.Range("G37").Value = 24
$ws.Range("H37").Value = 657

$ws.Range("A38").Value = "Austria"
$ws.Range("B38").Value = 15961
$ws.Range("C38").Value = 79
$ws.Range("D38").Value = 14148
$ws.Range("E38").Value = 1190
$ws.Range("F38").Value = 59
$ws.Range("G38").Value = 3
$ws.Range("H38").Value = 623

# Row 60: Kazajistan - refreshed totals
$ws.Range("D60").Value = 2223
$ws.Range("E60").Value = 3024

# Row 100: El Salvador - refreshed totals
$ws.Range("E100").Value = 629
$ws.Range("G100").Value = 2
$ws.Range("H100").Value = 20

# Rows 109-111: Guinea-Bisau moves above Costa Rica (table sorted by total cases desc)
$ws.Range("A109").Value = "Guinea-Bisau"
$ws.Range("B109").Value = 820
$ws.Range("C109").Value = 59
$ws.Range("D109").Value = 26
$ws.Range("E109").Value = 791
$ws.Range("F109").Value = 0
$ws.Range("G109").Value = 0
$ws.Range("H109").Value = 3

$ws.Range("A110").Value = "Costa Rica"
$ws.Range("B110").Value = 804
$ws.Range("C110").Value = 3
$ws.Range("D110").Value = 520
$ws.Range("E110").Value = 277
$ws.Range("F110").Value = 6
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 7

# Row 111: Burkina Faso - refreshed totals (stays in place)
$ws.Range("B111").Value = 766
$ws.Range("C111").Value = 6
$ws.Range("D111").Value = 588
$ws.Range("E111").Value = 127
$ws.Range("G111").Value = 1
$ws.Range("H111").Value = 51

# Row 167: Islas Caimanes - refreshed totals
$ws.Range("B167").Value = 85
$ws.Range("C167").Value = 1
$ws.Range("D167").Value = 50
$ws.Range("E167").Value = 34

# Row 168: Barbados - refreshed totals
$ws.Range("B168").Value = 85
$ws.Range("C168").Value = 1
$ws.Range("E168").Value = 21

# Row 196: San Vicente y las Granadinas - refreshed totals
$ws.Range("D196").Value = 12
$ws.Range("E196").Value = 5
